$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.498.38"
$ws.Range("E2").Value = "  -0.04%  "

$ws.Range("D3").Value = "3.580.66"
$ws.Range("E3").Value = "  +0.64%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'604.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.01%  "

$ws.Range("D6").Value = "'135.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.51%  "

$ws.Range("D7").Value = "3.579.96"
$ws.Range("E7").Value = "  +0.62%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").Value = "'0.495"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.75%  "

$ws.Range("D10").Value = "'0.124"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.91%  "

$ws.Range("D11").Value = "'7.19"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.56%  "

$ws.Range("D12").Value = "'0.391"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.22%  "

$ws.Range("D13").Value = "4.197.05"
$ws.Range("E13").Value = "  +0.83%  "

$ws.Range("D14").Value = "'0.0000185"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.45%  "

$ws.Range("D15").Value = "'27.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.55%  "

$ws.Range("D16").Value = "3.585.14"
$ws.Range("E16").Value = "  +0.60%  "

$ws.Range("D17").Value = "'0.116"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.05%  "

$ws.Range("D18").Value = "65.607.21"
$ws.Range("E18").Value = "  +0.10%  "

$ws.Range("D19").Value = "'10.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.99%  "

$ws.Range("D20").Value = "'14.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.41%  "

$ws.Range("E21").Value = "  -0.51%  "

$ws.Range("D22").Value = "'394.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.26%  "

$ws.Range("D23").Value = "'0.585"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.12%  "

$ws.Range("D24").Value = "3.726.55"
$ws.Range("E24").Value = "  +0.75%  "

$ws.Range("D25").Value = "'74.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.25%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("E27").Value = "  -0.93%  "

$ws.Range("D28").Value = "'8.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.99%  "

$ws.Range("E29").Value = "  +27.21%  "

$ws.Range("E30").Value = "  +3.49%  "

$ws.Range("D31").Value = "'8.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.34%  "

$ws.Range("D32").Value = "'1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.00%  "

$ws.Range("D33").Value = "3.579.89"
$ws.Range("E33").Value = "  +0.25%  "

$ws.Range("D34").Value = "'24.39"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.37%  "

$ws.Range("E35").Value = "  +0.44%  "

$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("D37").Value = "'5.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.01%  "

$ws.Range("D38").Value = "'1.60"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.43%  "

$ws.Range("D39").Value = "'7.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.10%  "

$ws.Range("D40").Value = "'170.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.92%  "

$ws.Range("D41").Value = "'0.0828"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.92%  "

$ws.Range("D42").Value = "'0.835"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.09%  "

$ws.Range("D43").Value = "'26.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.42%  "

$ws.Range("D44").Value = "'43.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.44%  "

$ws.Range("D45").Value = "'1.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.19%  "

$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").Value = "'1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.11%  "

$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "'4.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.09%  "

$ws.Range("E48").Value = "  -1.08%  "

$ws.Range("D49").Value = "'7.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.65%  "

$ws.Range("D50").Value = "2.457.39"
$ws.Range("E50").Value = "  -1.26%  "

$ws.Range("E51").Value = "  +2.00%  "
